$wb = $excel.ActiveWorkbook

# Plain (non-numeric-looking) text: a normal .Value assignment is stored as a
# shared string by the engine without any style churn.
function Set-PlainText($range, $text) {
    $range.Value = $text
}

# Digit-only text (e.g. "760604", "062637") would otherwise be auto-coerced
# to a number by .Value. Forcing NumberFormat="@" first makes the engine
# store it as a real (shared) string; we then restore the original
# NumberFormat so the cell's style index is unchanged from before the edit.
# "0" round-trips onto the workbook's existing built-in "0" number format
# (style index 3 here); omitting the restore step (i.e. leaving NumberFormat
# "@") is what's needed for cells whose original style was the default style
# (index 0) - achieved via resetting to the "Normal" cell style.
function Set-NumericLookingText($range, $text, $restoreFormat) {
    $range.NumberFormat = "@"
    $range.Value = $text
    if ($restoreFormat -eq "Normal") {
        $range.Style = "Normal"
    } else {
        $range.NumberFormat = $restoreFormat
    }
}

# --- Customer sheet ---
$ws = $wb.Worksheets.Item("Customer")
Set-PlainText $ws.Range("A2") "test_brptge"
Set-NumericLookingText $ws.Range("B2") "545960" "Normal"
Set-PlainText $ws.Range("A4") "test_ceetfd"
Set-NumericLookingText $ws.Range("B4") "820723" "Normal"

# --- BA sheet ---
$ws = $wb.Worksheets.Item("BA")
Set-NumericLookingText $ws.Range("A2") "760604" "0"
Set-PlainText $ws.Range("B2") "test_mphoyo"
Set-NumericLookingText $ws.Range("A4") "760604" "0"
Set-PlainText $ws.Range("B4") "test_oovgak"

# --- User sheet ---
$ws = $wb.Worksheets.Item("User")
Set-PlainText $ws.Range("A2") "test_zpoyjz"
Set-NumericLookingText $ws.Range("B2") "634223" "Normal"
Set-PlainText $ws.Range("A3") "test_kvauhi"
Set-NumericLookingText $ws.Range("B3") "261002" "Normal"
Set-PlainText $ws.Range("A5") "test_wlmvuo"
Set-NumericLookingText $ws.Range("B5") "062637" "Normal"

# --- WithholdingTax sheet ---
$ws = $wb.Worksheets.Item("WithholdingTax")
Set-NumericLookingText $ws.Range("B2") "81" "Normal"
Set-NumericLookingText $ws.Range("C2") "47" "Normal"
